$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 60, shifting existing rows 60..116 down to 61..117.
$ws.Rows.Item(60).EntireRow.Insert()

# Fill in the new row 60 with the latest weekly price entry.
$ws.Cells.Item(60, 1).Value = 8
$ws.Cells.Item(60, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(60, 3).Value = "Coquimbo"
$ws.Cells.Item(60, 4).Value = 44447
$ws.Cells.Item(60, 5).Value = 4
$ws.Cells.Item(60, 6).Value = 100112031
$ws.Cells.Item(60, 7).Value = "Poroto verde"
$ws.Cells.Item(60, 8).Value = "Magnum"
$ws.Cells.Item(60, 9).Value = "Primera"
$ws.Cells.Item(60, 10).Value = 600
$ws.Cells.Item(60, 11).Value = 34000
$ws.Cells.Item(60, 12).Value = 35000
$ws.Cells.Item(60, 13).Value = 34500
$ws.Cells.Item(60, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(60, 15).Value = "Perú"
$ws.Cells.Item(60, 16).Value = 1380
$ws.Cells.Item(60, 17).Value = 25
$ws.Cells.Item(60, 18).Value = "Hortaliza"
